# ---------------------------------------------------------------------------
# Adds a second worksheet ("correl-employment-GES") that analyses the
# correlation between Singapore's overall employment and graduate
# employment, mirroring the existing "correl-gdp-employment" sheet, plus a
# new scatter chart, and extends the original sheet with a further data
# point (2019).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Extend "correl-gdp-employment" with the 2019 data point.
# ---------------------------------------------------------------------------
$ws1.Range("A8").Value = 2019
$ws1.Range("A8").NumberFormat = "0"

$ws1.Range("B8").Value = 510737.8
$ws1.Range("C8").Value = 3631.7

$ws1.Range("B10").Formula = "=CORREL(B2:B8,C2:C8)"

# ---------------------------------------------------------------------------
# 2. Add the new worksheet right after the first one.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "correl-employment-GES"

$ws2.Columns.Item(1).ColumnWidth = 16.15625
$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(3).ColumnWidth = 15.5234375

# Header row
$ws2.Range("A1").Value = "Year"

$ws2.Range("B1").Formula = "=[2]processed!`$C1"
$ws2.Range("C1").Value = "Employed Graduates"

$ws2.Range("A1:C1").WrapText = $true

# Year column
$ws2.Range("A2").Value = 2013
$ws2.Range("A3").Value = 2014
$ws2.Range("A4").Value = 2015
$ws2.Range("A5").Value = 2016
$ws2.Range("A6").Value = 2017
$ws2.Range("A7").Value = 2018
$ws2.Range("A8").Value = 2019

# Employed (thousands) - these cells resolve against the already-cached
# external link data, so they can stay live formulas.
$ws2.Range("B2").Formula = "=[2]processed!`$C2"
$ws2.Range("B3").Formula = "=[2]processed!`$C3"
$ws2.Range("B4").Formula = "=[2]processed!`$C4"
$ws2.Range("B5").Formula = "=[2]processed!`$C5"
$ws2.Range("B6").Formula = "=[2]processed!`$C6"
$ws2.Range("B7").Formula = "=[2]processed!`$C7"
$ws2.Range("B8").Value = 3631.7

$ws2.Range("B2:B8").WrapText = $true

# Employed Graduates
$ws2.Range("C2").Value = 10029
$ws2.Range("C3").Value = 10167
$ws2.Range("C4").Value = 10337
$ws2.Range("C5").Value = 10944
$ws2.Range("C6").Value = 12551
$ws2.Range("C7").Value = 12626
$ws2.Range("C8").Value = 12900

$ws2.Range("C2:C8").WrapText = $true

# Correlation summary row
$ws2.Range("A10").Value = "Correlation between Singapore Overall Employment and Graduate Employment"
$ws2.Range("A10").WrapText = $true

$ws2.Range("B10").Formula = "=CORREL(B2:B8,C2:C8)"
$ws2.Range("B10").HorizontalAlignment = -4108
$ws2.Range("B10").VerticalAlignment = -4108
$ws2.Range("B10").Interior.Color = 65535

$ws2.Range("C10").Value = "The correlation between these 2 factors are very strong"
$ws2.Range("C10").HorizontalAlignment = -4108
$ws2.Range("C10").VerticalAlignment = -4108
$ws2.Range("C10").WrapText = $true
$ws2.Range("C10").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. Chart: Overall VS Graduate Employment (Singapore)
# ---------------------------------------------------------------------------
$co = $ws2.ChartObjects().Add(209550, 245744, 4800000, 3500000)
$chart = $co.Chart
$chart.ChartType = 74

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Overall VS Graduate Employment (Singapore)"

$ser1 = $chart.SeriesCollection().NewSeries()
$ser1.Name = "='correl-employment-GES'!`$B`$1"
$ser1.XValues = $ws2.Range("A2:A8")
$ser1.Values = $ws2.Range("B2:B8")
$ser1.MarkerStyle = 8
$ser1.HasDataLabels = $true

$ser2 = $chart.SeriesCollection().NewSeries()
$ser2.Name = "='correl-employment-GES'!`$C`$1"
$ser2.XValues = $ws2.Range("A2:A8")
$ser2.Values = $ws2.Range("C2:C8")
$ser2.MarkerStyle = 8
$ser2.HasDataLabels = $true

$chart.HasLegend = $true
$chart.Legend.Position = -4107

$ws2.Range("A1").Select()
